$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "249.30"
Set-TextValue "G2" "11"
Set-TextValue "D3" "22.52"
Set-TextValue "G3" "11"
Set-TextValue "D4" "5.354"
Set-TextValue "G4" "11"
Set-TextValue "D5" "0.05689"
Set-TextValue "G5" "11"
Set-TextValue "D6" "3.400"
Set-TextValue "G6" "11"
Set-TextValue "D7" "6.316"
Set-TextValue "G7" "11"
Set-TextValue "D8" "0.8062"
Set-TextValue "G8" "11"
Set-TextValue "D9" "0.9148"
Set-TextValue "G9" "11"
Set-TextValue "D10" "0.1407"
Set-TextValue "G10" "11"
Set-TextValue "D11" "0.07445"
Set-TextValue "G11" "11"
Set-TextValue "D12" "0.03106"
Set-TextValue "G12" "11"
Set-TextValue "D13" "0.03007"
Set-TextValue "G13" "11"
Set-TextValue "D14" "0.09388"
Set-TextValue "G14" "11"
Set-TextValue "D15" "3.868"
Set-TextValue "G15" "11"
Set-TextValue "D16" "0.001591"
Set-TextValue "G16" "11"
Set-TextValue "D17" "0.04799"
Set-TextValue "G17" "11"
Set-TextValue "G18" "11"
Set-TextValue "D19" "0.0005851"
Set-TextValue "G19" "11"
Set-TextValue "D20" "0.006468"
Set-TextValue "G20" "11"
Set-TextValue "D21" "0.004991"
Set-TextValue "G21" "11"
Set-TextValue "D22" "0.0009916"
Set-TextValue "G22" "11"
Set-TextValue "D23" "0.0001500"
Set-TextValue "G23" "11"
Set-TextValue "D24" "3.700"
Set-TextValue "G24" "11"
Set-TextValue "D25" "2.197"
Set-TextValue "G25" "11"
Set-TextValue "D26" "0.3254"
Set-TextValue "G26" "11"
Set-TextValue "G27" "11"
Set-TextValue "G28" "11"
Set-TextValue "G29" "11"
Set-TextValue "G30" "11"
Set-TextValue "G31" "11"
Set-TextValue "G32" "11"
Set-TextValue "G33" "11"
Set-TextValue "G34" "11"
Set-TextValue "G35" "11"
Set-TextValue "G36" "11"
Set-TextValue "G37" "11"
Set-TextValue "G38" "11"
Set-TextValue "G39" "11"
Set-TextValue "D40" "0.04002"
Set-TextValue "G40" "11"
Set-TextValue "B41" "BKEXToken"
Set-TextValue "C41" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1072"
Set-TextValue "E41" "40BKEXTokenBKK"
Set-TextValue "G41" "11"
Set-TextValue "B42" "CEJI"
Set-TextValue "C42" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.002730"
Set-TextValue "E42" "41CEJICEJI"
Set-TextValue "G42" "11"
Set-TextValue "B43" "KickToken"
Set-TextValue "C43" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.003035"
Set-TextValue "E43" "42KickTokenKICKWorstin24h"
Set-TextValue "G43" "11"
Set-TextValue "D44" "0.007961"
Set-TextValue "G44" "11"
Set-TextValue "D45" "0.00005752"
Set-TextValue "G45" "11"
Set-TextValue "G46" "11"
Set-TextValue "D47" "0.4991"
Set-TextValue "E47" "46CoinbaseStockTokenCOIN"
Set-TextValue "G47" "11"
Set-TextValue "D48" "0.2067"
Set-TextValue "G48" "11"
Set-TextValue "D49" "0.00002100"
Set-TextValue "G49" "11"
Set-TextValue "G50" "11"
Set-TextValue "G51" "11"
